# Daily attendance processing - 2025-11-10 21:21:42
#
# Column G ("Recorded By") lists the users who recorded each attendance
# session, as a comma-separated string (e.g. "dnasr281@gmail.com, System").
# Going forward "System, system" entries should be listed first in that
# column, with the remaining recorder(s) following in their original
# relative order.
#
# Walk every data row in the used range and, wherever "System"/"system"
# is present but not already first, move it (them) to the front while
# keeping the relative order of everything else unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Column G is "Recorded By" (column index 7); row 1 is the header.
$col = 7

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -gt 1) {
            $systemParts = @()
            $otherParts = @()

            foreach ($p in $parts) {
                if ($p.ToLower() -eq "system") {
                    $systemParts += $p
                } else {
                    $otherParts += $p
                }
            }

            if ($systemParts.Length -gt 0) {
                $newParts = $systemParts + $otherParts
                $newVal = $newParts -join ", "

                if ($newVal -ne $val) {
                    $cell.Value2 = $newVal
                }
            }
        }
    }
}
